# Update "想去人数" (interested count) values per the diff for gh-pages output at 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2744
$ws.Range("F4").Value = 1078
$ws.Range("F5").Value = 19901
$ws.Range("F7").Value = 2309
$ws.Range("F8").Value = 755
$ws.Range("F10").Value = 444
$ws.Range("F11").Value = 697
$ws.Range("F12").Value = 245
$ws.Range("F15").Value = 379
$ws.Range("F16").Value = 82
$ws.Range("F17").Value = 269
$ws.Range("F18").Value = 172
$ws.Range("F19").Value = 205

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 19
$ws.Range("F8").Value = 135
$ws.Range("F11").Value = 11
$ws.Range("F14").Value = 90
$ws.Range("F16").Value = 94

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6013
$ws.Range("F3").Value = 652
$ws.Range("F4").Value = 596

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6013
$ws.Range("F3").Value = 652
$ws.Range("F4").Value = 596
$ws.Range("F8").Value = 2744
$ws.Range("F9").Value = 1078
$ws.Range("F10").Value = 19901
$ws.Range("F12").Value = 19
$ws.Range("F16").Value = 2309
$ws.Range("F17").Value = 755
$ws.Range("F18").Value = 135
$ws.Range("F20").Value = 444
$ws.Range("F21").Value = 697
$ws.Range("F22").Value = 245
$ws.Range("F27").Value = 11
$ws.Range("F28").Value = 379
$ws.Range("F29").Value = 82
$ws.Range("F32").Value = 269
$ws.Range("F33").Value = 90
$ws.Range("F34").Value = 172
$ws.Range("F36").Value = 205
$ws.Range("F37").Value = 94
$ws.Range("F38").Value = 94
